$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7161694356617545
$ws.Range("C2").Value = 0.2154196615050452
$ws.Range("E2").Value = 0.2232326671488742
$ws.Range("F2").Value = 2.295485513928426
$ws.Range("G2").Value = 0.002495197373909955
$ws.Range("J2").Value = 0.08280919888975191
$ws.Range("K2").Value = 0.3006274835378235
$ws.Range("L2").Value = 0.3885996225053816
$ws.Range("M2").Value = 0.254535608282378
$ws.Range("N2").Value = 2.310740456998725
$ws.Range("O2").Value = 4.436040179096722
$ws.Range("B3").Value = 0.6829294984964349
$ws.Range("C3").Value = 0.2157782235239303
$ws.Range("E3").Value = 0.2231367459586835
$ws.Range("F3").Value = 2.297476359903754
$ws.Range("G3").Value = 0.002497500711587901
$ws.Range("J3").Value = 0.08140212343493403
$ws.Range("K3").Value = 0.2723230496631714
$ws.Range("L3").Value = 0.3840320389043086
$ws.Range("M3").Value = 0.2472488178745955
$ws.Range("N3").Value = 2.331335495937918
$ws.Range("O3").Value = 4.46054356098702
$ws.Range("B4").Value = 0.662791503486261
$ws.Range("C4").Value = 0.2160165652860186
$ws.Range("E4").Value = 0.223149177781746
$ws.Range("F4").Value = 2.299689748693119
$ws.Range("G4").Value = 0.002498991340806535
$ws.Range("J4").Value = 0.08052860848033205
$ws.Range("K4").Value = 0.2549974999389946
$ws.Range("L4").Value = 0.3813725297014869
$ws.Range("M4").Value = 0.2428759042067554
$ws.Range("N4").Value = 2.344623970058002
$ws.Range("O4").Value = 4.477471367688025
$ws.Range("B5").Value = 0.6546540067851936
$ws.Range("C5").Value = 0.2161182789435436
$ws.Range("E5").Value = 0.2231722330739103
$ws.Range("F5").Value = 2.300841291247941
$ws.Range("G5").Value = 0.002499618044358443
$ws.Range("J5").Value = 0.0801702523767851
$ws.Range("K5").Value = 0.2479510638619047
$ws.Range("L5").Value = 0.3803253800952788
$ws.Range("M5").Value = 0.241119522237824
$ws.Range("N5").Value = 2.350200892119759
$ws.Range("O5").Value = 4.484843345705542
$ws.Range("B6").Value = 0.6533069618307934
$ws.Range("C6").Value = 0.2161354459363221
$ws.Range("E6").Value = 0.2231771494582411
$ws.Range("F6").Value = 2.301047588326355
$ws.Range("G6").Value = 0.002499723272754977
$ws.Range("J6").Value = 0.08011060356383837
$ws.Range("K6").Value = 0.2467818582559005
$ws.Range("L6").Value = 0.3801537182584624
$ws.Range("M6").Value = 0.2408294283975501
$ws.Range("N6").Value = 2.351136707666297
$ws.Range("O6").Value = 4.48609608145712
$ws.Range("B7").Value = 0.6626814785406623
$ws.Range("C7").Value = 0.2160179184370143
$ws.Range("E7").Value = 0.2231494158048832
$ws.Range("F7").Value = 2.299704267809467
$ws.Range("G7").Value = 0.00249899971453986
$ws.Range("J7").Value = 0.08052378523061776
$ws.Range("K7").Value = 0.25490241254262
$ws.Range("L7").Value = 0.381358258982992
$ws.Range("M7").Value = 0.242852113086041
$ws.Range("N7").Value = 2.344698527330427
$ws.Range("O7").Value = 4.477568870057453
$ws.Range("B8").Value = 0.7046523352832708
$ws.Range("C8").Value = 0.2155395288111421
$ws.Range("E8").Value = 0.2231848209647787
$ws.Range("F8").Value = 2.295966494658316
$ws.Range("G8").Value = 0.002495975747738094
$ws.Range("J8").Value = 0.0823260322419479
$ws.Range("K8").Value = 0.2908572543815211
$ws.Range("L8").Value = 0.3869947171085073
$ws.Range("M8").Value = 0.2520022152831736
$ws.Range("N8").Value = 2.317708204541884
$ws.Range("O8").Value = 4.444098415863294
$ws.Range("B9").Value = 0.7890875301921483
$ws.Range("C9").Value = 0.214745020757185
$ws.Range("E9").Value = 0.2238180374076428
$ws.Range("F9").Value = 2.296485555555833
$ws.Range("G9").Value = 0.002490649169345293
$ws.Range("J9").Value = 0.0857838874231085
$ws.Range("K9").Value = 0.3617742255355267
$ws.Range("L9").Value = 0.3991922295393522
$ws.Range("M9").Value = 0.2707423210155184
$ws.Range("N9").Value = 2.269878801944893
$ws.Range("O9").Value = 4.393386703088623
$ws.Range("B10").Value = 0.8523961493472996
$ws.Range("C10").Value = 0.214247956138685
$ws.Range("E10").Value = 0.224624302377844
$ws.Range("F10").Value = 2.301635564790161
$ws.Range("G10").Value = 0.002487100007837948
$ws.Range("J10").Value = 0.08827743201822358
$ws.Range("K10").Value = 0.4141120157540854
$ws.Range("L10").Value = 0.4088444819549437
$ws.Range("M10").Value = 0.2849896644871848
$ws.Range("N10").Value = 2.237840546420596
$ws.Range("O10").Value = 4.365210116362334
$ws.Range("B11").Value = 0.8814686961084988
$ws.Range("C11").Value = 0.214040459741728
$ws.Range("E11").Value = 0.2250645999727681
$ws.Range("F11").Value = 2.305010677944054
$ws.Range("G11").Value = 0.002485563740316943
$ws.Range("J11").Value = 0.08940154230024433
$ws.Range("K11").Value = 0.4379700807352265
$ws.Range("L11").Value = 0.4133841346223619
$ws.Range("M11").Value = 0.2915738218073614
$ws.Range("N11").Value = 2.223937762736311
$ws.Range("O11").Value = 4.35436048395394
$ws.Range("B12").Value = 0.8925164300922006
$ws.Range("C12").Value = 0.213964549048054
$ws.Range("E12").Value = 0.225241851592866
$ws.Range("F12").Value = 2.306436871421397
$ws.Range("G12").Value = 0.002484993192534023
$ws.Range("J12").Value = 0.08982573287253359
$ws.Range("K12").Value = 0.4470112483433297
$ws.Range("L12").Value = 0.4151244350272236
$ws.Range("M12").Value = 0.2940817294224303
$ws.Range("N12").Value = 2.2187696527887
$ws.Range("O12").Value = 4.350534721407314
$ws.Range("B13").Value = 0.8901353974623873
$ws.Range("C13").Value = 0.2139807795284021
$ws.Range("E13").Value = 0.2252032102138273
$ws.Range("F13").Value = 2.306123133249088
$ws.Range("G13").Value = 0.002485115572684655
$ws.Range("J13").Value = 0.08973444218094073
$ws.Range("K13").Value = 0.4450637832991902
$ws.Range("L13").Value = 0.4147486884909597
$ws.Range("M13").Value = 0.2935409587819109
$ws.Range("N13").Value = 2.219878402613258
$ws.Range("O13").Value = 4.351346095962413
$ws.Range("B14").Value = 0.8823768303298323
$ws.Range("C14").Value = 0.2140341612183647
$ws.Range("E14").Value = 0.2250789720226223
$ws.Range("F14").Value = 2.305125045310959
$ws.Range("G14").Value = 0.002485516576760848
$ws.Range("J14").Value = 0.08943647058277548
$ws.Range("K14").Value = 0.4387137726672847
$ws.Range("L14").Value = 0.4135268855939387
$ws.Range("M14").Value = 0.2917798566793763
$ws.Range("N14").Value = 2.223510643752456
$ws.Range("O14").Value = 4.354040070671033
$ws.Range("B15").Value = 0.8776294893353054
$ws.Range("C15").Value = 0.2140672055261241
$ws.Range("E15").Value = 0.2250042410732505
$ws.Range("F15").Value = 2.304532966032369
$ws.Range("G15").Value = 0.002485763660690623
$ws.Range("J15").Value = 0.08925376041212019
$ws.Range("K15").Value = 0.4348250583984168
$ws.Range("L15").Value = 0.4127812562636564
$ws.Range("M15").Value = 0.2907030308182996
$ws.Range("N15").Value = 2.225748072974726
$ws.Range("O15").Value = 4.355727024165219
$ws.Range("B16").Value = 0.8505016276993445
$ws.Range("C16").Value = 0.2142618900176103
$ws.Range("E16").Value = 0.2245970029208024
$ws.Range("F16").Value = 2.301435737721505
$ws.Range("G16").Value = 0.002487201977143271
$ws.Range("J16").Value = 0.08820376178675104
$ws.Range("K16").Value = 0.4125537887441908
$ws.Range("L16").Value = 0.40855078508757
$ws.Range("M16").Value = 0.2845614312236648
$ws.Range("N16").Value = 2.238762639876502
$ws.Range("O16").Value = 4.365958746452606
$ws.Range("B17").Value = 0.8339290329462017
$ws.Range("C17").Value = 0.2143860819019707
$ws.Range("E17").Value = 0.2243659695635749
$ws.Range("F17").Value = 2.299799856519428
$ws.Range("G17").Value = 0.002488104346544118
$ws.Range("J17").Value = 0.08755699339965162
$ws.Range("K17").Value = 0.3989033972424068
$ws.Range("L17").Value = 0.4059935212580257
$ws.Range("M17").Value = 0.2808200060098116
$ws.Range("N17").Value = 2.246918661959447
$ws.Range("O17").Value = 4.372739464646429
$ws.Range("B18").Value = 0.8244226648516246
$ws.Range("C18").Value = 0.2144592672952363
$ws.Range("E18").Value = 0.2242400087927834
$ws.Range("F18").Value = 2.298956130706884
$ws.Range("G18").Value = 0.002488630735244792
$ws.Range("J18").Value = 0.08718402878366049
$ws.Range("K18").Value = 0.3910567248305767
$ws.Range("L18").Value = 0.4045366658116194
$ws.Range("M18").Value = 0.2786777397172244
$ws.Range("N18").Value = 2.251673002175687
$ws.Range("O18").Value = 4.376824809952268
$ws.Range("B19").Value = 0.8212084154384343
$ws.Range("C19").Value = 0.214484348204202
$ws.Range("E19").Value = 0.2241985513477722
$ws.Range("F19").Value = 2.298687162407873
$ws.Range("G19").Value = 0.002488810228686887
$ws.Range("J19").Value = 0.08705758490177118
$ws.Range("K19").Value = 0.3884007927351831
$ws.Range("L19").Value = 0.4040458116867001
$ws.Range("M19").Value = 0.2779540769618762
$ws.Range("N19").Value = 2.253293596854338
$ws.Range("O19").Value = 4.378239862949812
$ws.Range("B20").Value = 0.8356905538776971
$ws.Range("C20").Value = 0.2143726800855887
$ws.Range("E20").Value = 0.2243898472799088
$ws.Range("F20").Value = 2.299963942541908
$ws.Range("G20").Value = 0.002488007525424757
$ws.Range("J20").Value = 0.08762594253096978
$ws.Range("K20").Value = 0.4003560245375866
$ws.Range("L20").Value = 0.4062642970394421
$ws.Range("M20").Value = 0.2812172842150318
$ws.Range("N20").Value = 2.246043897282564
$ws.Range("O20").Value = 4.371998474131345
$ws.Range("B21").Value = 0.8846546676707305
$ws.Range("C21").Value = 0.2140184095414277
$ws.Range("E21").Value = 0.2251151786835557
$ws.Range("F21").Value = 2.305414190936716
$ws.Range("G21").Value = 0.002485398488588063
$ws.Range("J21").Value = 0.08952403251564789
$ws.Range("K21").Value = 0.4405787482731682
$ws.Range("L21").Value = 0.4138851837868316
$ws.Range("M21").Value = 0.2922967394082363
$ws.Range("N21").Value = 2.222441145517303
$ws.Range("O21").Value = 4.353241112967567
$ws.Range("B22").Value = 0.9168801688197448
$ws.Range("C22").Value = 0.2138023928907273
$ws.Range("E22").Value = 0.2256505219080474
$ws.Range("F22").Value = 2.309839392093465
$ws.Range("G22").Value = 0.002483758610910651
$ws.Range("J22").Value = 0.09075587706206534
$ws.Range("K22").Value = 0.4669050908795214
$ws.Range("L22").Value = 0.4189895566480573
$ws.Range("M22").Value = 0.2996229973496938
$ws.Range("N22").Value = 2.207578283397583
$ws.Range("O22").Value = 4.342630109654294
$ws.Range("B23").Value = 0.8996604860245725
$ws.Range("C23").Value = 0.213916269572362
$ws.Range("E23").Value = 0.2253592076687987
$ws.Range("F23").Value = 2.307398712929711
$ws.Range("G23").Value = 0.002484627888115797
$ws.Range("J23").Value = 0.09009921699318113
$ws.Range("K23").Value = 0.4528508583851192
$ws.Range("L23").Value = 0.4162539940320613
$ws.Range("M23").Value = 0.2957050988272272
$ws.Range("N23").Value = 2.215459366379451
$ws.Range("O23").Value = 4.348142688397076
$ws.Range("B24").Value = 0.8348941034376196
$ws.Range("C24").Value = 0.2143787334839971
$ws.Range("E24").Value = 0.2243790307820497
$ws.Range("F24").Value = 2.299889457813578
$ws.Range("G24").Value = 0.002488051274776234
$ws.Range("J24").Value = 0.08759477414353611
$ws.Range("K24").Value = 0.3996992882159986
$ws.Range("L24").Value = 0.4061418377142871
$ws.Range("M24").Value = 0.2810376475357899
$ws.Range("N24").Value = 2.246439174877429
$ws.Range("O24").Value = 4.372332893365922
$ws.Range("B25").Value = 0.7660200447927537
$ws.Range("C25").Value = 0.2149446739977101
$ws.Range("E25").Value = 0.223586653333161
$ws.Range("F25").Value = 2.295506539147922
$ws.Range("G25").Value = 0.002492025922393029
$ws.Range("J25").Value = 0.08485666007418047
$ws.Range("K25").Value = 0.342546963326555
$ws.Range("L25").Value = 0.3957706764585254
$ws.Range("M25").Value = 0.2655880224525689
$ws.Range("N25").Value = 2.282272741409731
$ws.Range("O25").Value = 4.405509599261364
